$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" '35.213.23'
Set-TextValue "E2" '  +0.20%  '
Set-TextValue "D3" '1.890.65'
Set-TextValue "E3" '  +2.06%  '
Set-TextValue "E4" '  -0.25%  '
Set-TextValue "D5" '242.96'
Set-TextValue "E5" '  +2.24%  '
Set-TextValue "D6" '0.655'
Set-TextValue "E6" '  +5.56%  '
Set-TextValue "D8" '41.17'
Set-TextValue "E8" '  -1.96%  '
Set-TextValue "D9" '0.345'
Set-TextValue "E9" '  +5.65%  '
Set-TextValue "D10" '50.11'
Set-TextValue "E10" '  +7.82%  '
Set-TextValue "D11" '0.0708'
Set-TextValue "E11" '  +2.29%  '
Set-TextValue "D12" '0.0992'
Set-TextValue "E12" '  +0.48%  '
Set-TextValue "D13" '2.164.97'
Set-TextValue "E13" '  +2.04%  '
Set-TextValue "D14" '11.88'
Set-TextValue "E14" '  +4.47%  '
Set-TextValue "B15" 'Polygon'
Set-TextValue "C15" 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue "D15" '0.690'
Set-TextValue "E15" '  +2.25%  '
Set-TextValue "B16" 'WrappedEther'
Set-TextValue "C16" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue "D16" '1.874.61'
Set-TextValue "E16" '  +1.13%  '
Set-TextValue "E17" '  +2.35%  '
Set-TextValue "D18" '35.203.07'
Set-TextValue "E18" '  +0.29%  '
Set-TextValue "D19" '71.06'
Set-TextValue "D20" '0.0₃0812'
Set-TextValue "E20" '  +2.52%  '
Set-TextValue "D21" '240.13'
Set-TextValue "E21" '  -0.21%  '
Set-TextValue "D22" '12.38'
Set-TextValue "E22" '  +1.32%  '
Set-TextValue "E23" '  +0.19%  '
Set-TextValue "E24" '  -0.33%  '
Set-TextValue "D25" '2.42'
Set-TextValue "E25" '  +32.40%  '
Set-TextValue "E26" '  +0.62%  '
Set-TextValue "D27" '170.13'
Set-TextValue "E27" '  +1.05%  '
Set-TextValue "E28" '  +5.29%  '
Set-TextValue "D29" '18.20'
Set-TextValue "E29" '  +3.44%  '
Set-TextValue "E30" '  +2.37%  '
Set-TextValue "D31" '4.11'
Set-TextValue "E31" '  +3.46%  '
Set-TextValue "D32" '0.950'
Set-TextValue "E32" '  +15.80%  '
Set-TextValue "B33" 'Hedera'
Set-TextValue "C33" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D33" '0.0558'
Set-TextValue "E33" '  +0.82%  '
Set-TextValue "B34" 'BinanceUSD'
Set-TextValue "C34" 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue "D34" '1.01'
Set-TextValue "E34" '  -0.29%  '
Set-TextValue "E35" '  +2.32%  '
Set-TextValue "E36" '  -2.64%  '
Set-TextValue "D37" '2.02'
Set-TextValue "E37" '  +0.83%  '
Set-TextValue "D38" '1.32'
Set-TextValue "E38" '  +1.25%  '
Set-TextValue "D39" '0.0209'
Set-TextValue "E39" '  +4.41%  '
Set-TextValue "D40" '1.08'
Set-TextValue "E40" '  +1.48%  '
Set-TextValue "E41" '  +15.29%  '
Set-TextValue "D42" '16.06'
Set-TextValue "E42" '  +8.45%  '
Set-TextValue "D43" '89.00'
Set-TextValue "E43" '  -1.08%  '
Set-TextValue "D44" '1.332.34'
Set-TextValue "E44" '  -0.69%  '
Set-TextValue "D45" '48.35'
Set-TextValue "E45" '  +39.67%  '
Set-TextValue "D46" '2.36'
Set-TextValue "E46" '  +2.69%  '
Set-TextValue "E47" '  -1.09%  '
Set-TextValue "E48" '  +1.27%  '
Set-TextValue "D49" '6.52'
Set-TextValue "E49" '  +0.97%  '
Set-TextValue "D50" '2.075.90'
Set-TextValue "E50" '  +1.97%  '
Set-TextValue "D51" '11.13'
Set-TextValue "E51" '  -12.74%  '
